$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ideas")

# Insert a new row at position 21 and set its value
$ws.Rows.Item(21).Insert()
$ws.Range("A21").Value = "Add bullets fired in the score"

# Make "Ideas" the active sheet/tab and set its selection
$ws.Activate()
$ws.Range("A22").Select()
